$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values to reflect re-run with 5,000 simulation runs instead of 10,000.
$ws.Range("C2").Value = 55
$ws.Range("E2").Value = 6

$ws.Range("C4").Value = 109
$ws.Range("D4").Value = 48
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 39
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 33
$ws.Range("I4").Value = 44

$ws.Range("D5").Value = 59
$ws.Range("E5").Value = 11
$ws.Range("F5").Value = 39
$ws.Range("G5").Value = 52
$ws.Range("I5").Value = 54

$ws.Range("D6").Value = -92
$ws.Range("E6").Value = -18
$ws.Range("G6").Value = -81
$ws.Range("H6").Value = -39

$ws.Range("C7").Value = 51
$ws.Range("D7").Value = 98
$ws.Range("E7").Value = 19
$ws.Range("F7").Value = 64
$ws.Range("G7").Value = 86
$ws.Range("H7").Value = 41
$ws.Range("I7").Value = 90

$ws.Range("F8").Value = 12
$ws.Range("H8").Value = 3

$ws.Range("E9").Value = -30
$ws.Range("F9").Value = -31
$ws.Range("H9").Value = -2

$ws.Range("C10").Value = 13
$ws.Range("D10").Value = 19
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 6
$ws.Range("G10").Value = 17
$ws.Range("H10").Value = 8
$ws.Range("I10").Value = 18
